$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(65, 1).Value = "2024-10-04 00:00:00"
$ws.Cells.Item(65, 2).Value = 75650
$ws.Cells.Item(65, 3).Value = 10756.89
$ws.Cells.Item(65, 4).Value = 9519.370000000001
$ws.Cells.Item(65, 5).Value = 7.0494
